# aggiornamento fino a 28/06 incluso
# Appends rows 270-301 (dates 2021-05-28 .. 2021-06-28) to Sheet1,
# extending the dimension from A1:D269 to A1:D301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date-serial (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(270, 44344, 0, 0, 0),
    @(271, 44345, 1, 1, 30.53435114503817),
    @(272, 44346, 0, 1, 30.53435114503817),
    @(273, 44347, 1, 2, 61.06870229007634),
    @(274, 44348, 0, 2, 61.06870229007634),
    @(275, 44349, 0, 2, 61.06870229007634),
    @(276, 44350, 0, 2, 61.06870229007634),
    @(277, 44351, 0, 2, 61.06870229007634),
    @(278, 44352, 0, 1, 30.53435114503817),
    @(279, 44353, 0, 1, 30.53435114503817),
    @(280, 44354, 0, 0, 0),
    @(281, 44355, 0, 0, 0),
    @(282, 44356, 0, 0, 0),
    @(283, 44357, 0, 0, 0),
    @(284, 44358, 0, 0, 0),
    @(285, 44359, 0, 0, 0),
    @(286, 44360, 0, 0, 0),
    @(287, 44361, 0, 0, 0),
    @(288, 44362, 0, 0, 0),
    @(289, 44363, 0, 0, 0),
    @(290, 44364, 0, 0, 0),
    @(291, 44365, 0, 0, 0),
    @(292, 44366, 0, 0, 0),
    @(293, 44367, 0, 0, 0),
    @(294, 44368, 0, 0, 0),
    @(295, 44369, 0, 0, 0),
    @(296, 44370, 0, 0, 0),
    @(297, 44371, 0, 0, 0),
    @(298, 44372, 0, 0, 0),
    @(299, 44373, 0, 0, 0),
    @(300, 44374, 0, 0, 0),
    @(301, 44375, 0, 0, 0)
)

$lastExistingRow = 269

foreach ($entry in $data) {
    $r = $entry[0]

    # Copy the formatting (border/alignment/date number-format) of the
    # column-A cell from the last pre-existing row so the new date cells
    # keep the same style as every other row in the column.
    $ws.Cells.Item($lastExistingRow, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value2 = $entry[1]
    $ws.Cells.Item($r, 2).Value2 = $entry[2]
    $ws.Cells.Item($r, 3).Value2 = $entry[3]
    $ws.Cells.Item($r, 4).Value2 = $entry[4]
}

$excel.CutCopyMode = $false
